# Add "Test Executable ... - BootstrapFewShot" columns (AZ:BB) and a new
# data row (15) to the evaluation log sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- New header cells AZ1:BB1 -------------------------------------------
# Copy the formatting of the existing last header cell (AY1, bold/centered/
# bordered style) onto the three new header cells, then fill in their text.
$ws.Range("AY1").Copy()
$ws.Range("AZ1:BB1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("AZ1").Value = "Test Executable Time - BootstrapFewShot"
$ws.Range("BA1").Value = "Test Executable Scores - BootstrapFewShot"
$ws.Range("BB1").Value = "Test Executable Results - BootstrapFewShot"

# --- New data row 15 ------------------------------------------------------
$ws.Range("A15").Value = "mistral:7b-instruct-v0.3-q5_K_M"
$ws.Range("B15").Value = "llama3:70b"
$ws.Range("C15").Value = 1
$ws.Range("D15").Value = 200
$ws.Range("E15").Value = 1776.53

$ws.Range("AI15").Value = 1011.14

$ws.Range("AQ15").Value = 765.39
$ws.Range("AR15").Value = 71.25
$ws.Range("AS15").Value = "logs\mistral_7b_instruct_v0.3_q5_K_M_llama3_70b_1_200_test_bootstrap_match_2.txt"
$ws.Range("AT15").Value = 765.39
$ws.Range("AU15").Value = 42.5
$ws.Range("AV15").Value = "logs\mistral_7b_instruct_v0.3_q5_K_M_llama3_70b_1_200_test_bootstrap_correct_2.txt"
$ws.Range("AW15").Value = 54.57142857142857
$ws.Range("AX15").Value = 2
$ws.Range("AY15").Value = 2

$ws.Range("AZ15").Value = 765.39
$ws.Range("BA15").Value = 98.75
$ws.Range("BB15").Value = "logs\mistral_7b_instruct_v0.3_q5_K_M_llama3_70b_1_200_test_bootstrap_executable.txt"
